# Apply corrected "mitigation direct" results to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2" = 63.38
    "D2" = 32.77
    "Q2" = 0.15
    "U2" = 0.15

    "C3" = 64
    "D3" = 57

    "C4" = 65.15000000000001
    "D4" = 6.2
    "F4" = 2.01
    "M4" = 0.18
    "P4" = 0.18

    "C5" = 79
    "D5" = 50.25
    "F5" = 24.75
    "G5" = 0.75
    "K5" = 0.5
    "M5" = 0.5
    "N5" = 0.25
    "S5" = 0.25

    "C6" = 54.87
    "D6" = 64.87
    "F6" = 23.16
    "G6" = 0.26
    "K6" = 0.13
    "P6" = 0.26
    "U6" = 0.13
    "V6" = 0.66

    "C7" = 72.39
    "D7" = 2.1

    "C8" = 66.3
    "D8" = 44.81
    "F8" = 10.74
    "K8" = 0.37
    "N8" = 0.37
    "Q8" = 0.74

    "C9" = 53.94
    "D9" = 46.16
    "F9" = 28.69
    "G9" = 0.07000000000000001
    "I9" = 0.07000000000000001
    "R9" = 0.04

    "C10" = 73.43000000000001
    "D10" = 12.36
    "F10" = 18.46
    "G10" = 0.16
    "I10" = 0.08
    "L10" = 0.08
    "P10" = 0.32

    "C11" = 57.04
    "D11" = 47.04
    "K11" = 0.37
    "N11" = 0.37
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
